$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 57.365
$ws.Range("D2").Value = 57.365
$ws.Range("E2").Value = 2.97113999
$ws.Range("F2").Value = 0.02490991
$ws.Range("G2").Value = 1.41937062
$ws.Range("H2").Value = 81.53381680999999
$ws.Range("I2").Value = 5.320157805843646
$ws.Range("J2").Value = 5.320157805843646
$ws.Range("K2").Value = 0.2731317790016397
$ws.Range("L2").Value = 0.003154962008832622
$ws.Range("M2").Value = 0.1499306754151599
$ws.Range("N2").Value = 11.88823192923915

$ws.Range("C3").Value = 87.864
$ws.Range("D3").Value = 87.864
$ws.Range("E3").Value = 1.95999567
$ws.Range("F3").Value = 0.0166586
$ws.Range("G3").Value = 1.44007047
$ws.Range("H3").Value = 126.88050474
$ws.Range("I3").Value = 12.40763690045181
$ws.Range("J3").Value = 12.40763690045181
$ws.Range("K3").Value = 0.2644305860977794
$ws.Range("L3").Value = 0.002444669844917847
$ws.Range("M3").Value = 0.1277562330603281
$ws.Range("N3").Value = 23.3858671666646

$ws.Range("C4").Value = 32.682
$ws.Range("D4").Value = 65.304
$ws.Range("E4").Value = 2.64202544
$ws.Range("F4").Value = 0.04027917
$ws.Range("G4").Value = 0.65326478
$ws.Range("H4").Value = 21.63867518
$ws.Range("I4").Value = 4.834071187337257
$ws.Range("J4").Value = 9.640981927818999
$ws.Range("K4").Value = 0.3742397804925467
$ws.Range("L4").Value = 0.00591641203284917
$ws.Range("M4").Value = 0.1116700820881994
$ws.Range("N4").Value = 6.199053811861268

$ws.Range("C5").Value = 48.699
$ws.Range("D5").Value = 95.02200000000001
$ws.Range("E5").Value = 1.82448447
$ws.Range("F5").Value = 0.02652993
$ws.Range("G5").Value = 0.6324907099999999
$ws.Range("H5").Value = 31.05493616
$ws.Range("I5").Value = 8.541273461981671
$ws.Range("J5").Value = 15.39670661455152
$ws.Range("K5").Value = 0.2910290981907355
$ws.Range("L5").Value = 0.004417988605510252
$ws.Range("M5").Value = 0.07846775317341502
$ws.Range("N5").Value = 7.927988941195578

$ws.Range("C6").Value = 19.452
$ws.Range("D6").Value = 77.65000000000001
$ws.Range("E6").Value = 2.25203566
$ws.Range("F6").Value = 0.05408004
$ws.Range("G6").Value = 0.2606735
$ws.Range("H6").Value = 5.20954745
$ws.Range("I6").Value = 3.699375809698358
$ws.Range("J6").Value = 14.76128027629702
$ws.Range("K6").Value = 0.412703093308726
$ws.Range("L6").Value = 0.01178866740040949
$ws.Range("M6").Value = 0.06776512306972264
$ws.Range("N6").Value = 2.196665632224895

$ws.Range("C7").Value = 27.353
$ws.Range("D7").Value = 100.754
$ws.Range("E7").Value = 1.72111005
$ws.Range("F7").Value = 0.03848635
$ws.Range("G7").Value = 0.2579575
$ws.Range("H7").Value = 7.1603633
$ws.Range("I7").Value = 5.2057170966631
$ws.Range("J7").Value = 16.09162590990806
$ws.Range("K7").Value = 0.2817475500419696
$ws.Range("L7").Value = 0.007343530310069872
$ws.Range("M7").Value = 0.04594113851177132
$ws.Range("N7").Value = 2.271756111277273

$ws.Range("C8").Value = 13.249
$ws.Range("D8").Value = 79.294
$ws.Range("E8").Value = 2.23132336
$ws.Range("F8").Value = 0.06012335999999999
$ws.Range("G8").Value = 0.13249847
$ws.Range("H8").Value = 1.82867511
$ws.Range("I8").Value = 2.792991192216653
$ws.Range("J8").Value = 16.74520725359844
$ws.Range("K8").Value = 0.5006037627241369
$ws.Range("L8").Value = 0.01346878869012278
$ws.Range("M8").Value = 0.04056914367814383
$ws.Range("N8").Value = 0.8853916540604321

$ws.Range("C9").Value = 19.292
$ws.Range("D9").Value = 97.48
$ws.Range("E9").Value = 1.78270144
$ws.Range("F9").Value = 0.04292247
$ws.Range("G9").Value = 0.13589984
$ws.Range("H9").Value = 2.70927865
$ws.Range("I9").Value = 4.257565368008208
$ws.Range("J9").Value = 16.15916825578342
$ws.Range("K9").Value = 0.3042322848289751
$ws.Range("L9").Value = 0.008584551329022421
$ws.Range("M9").Value = 0.03343901599876057
$ws.Range("N9").Value = 1.270539684819084

$ws.Range("C10").Value = 9.757
$ws.Range("D10").Value = 77.81999999999999
$ws.Range("E10").Value = 2.309417
$ws.Range("F10").Value = 0.06146065
$ws.Range("G10").Value = 0.07538022000000001
$ws.Range("H10").Value = 0.78181485
$ws.Range("I10").Value = 2.365913892885485
$ws.Range("J10").Value = 18.81168465288209
$ws.Range("K10").Value = 0.6010924663799562
$ws.Range("L10").Value = 0.01479462070075401
$ws.Range("M10").Value = 0.0274361960885143
$ws.Range("N10").Value = 0.4406157677818572

$ws.Range("C11").Value = 14.788
$ws.Range("D11").Value = 91.646
$ws.Range("E11").Value = 1.91304846
$ws.Range("F11").Value = 0.04162814
$ws.Range("G11").Value = 0.07583028999999999
$ws.Range("H11").Value = 1.1796218
$ws.Range("I11").Value = 3.847838596129243
$ws.Range("J11").Value = 17.68274196394223
$ws.Range("K11").Value = 0.3704606071601755
$ws.Range("L11").Value = 0.00926349372841891
$ws.Range("M11").Value = 0.0227010566976758
$ws.Range("N11").Value = 0.6851812841761992

$ws.Range("C12").Value = 7.416
$ws.Range("D12").Value = 73.866
$ws.Range("E12").Value = 2.45605754
$ws.Range("F12").Value = 0.06082523000000001
$ws.Range("G12").Value = 0.04636355
$ws.Range("H12").Value = 0.37697081
$ws.Range("I12").Value = 1.949036362601464
$ws.Range("J12").Value = 19.36339290994698
$ws.Range("K12").Value = 0.6760375782667537
$ws.Range("L12").Value = 0.01633279295996165
$ws.Range("M12").Value = 0.02119164856414825
$ws.Range("N12").Value = 0.2650125554772915

$ws.Range("C13").Value = 12.024
$ws.Range("D13").Value = 84.048
$ws.Range("E13").Value = 2.0988823
$ws.Range("F13").Value = 0.04014076
$ws.Range("G13").Value = 0.04787667999999999
$ws.Range("H13").Value = 0.6141641999999999
$ws.Range("I13").Value = 3.342243946003133
$ws.Range("J13").Value = 17.69182570275597
$ws.Range("K13").Value = 0.4335243177721519
$ws.Range("L13").Value = 0.009542757098209901
$ws.Range("M13").Value = 0.01674002458183979
$ws.Range("N13").Value = 0.3894769235103758
